$wb = $excel.ActiveWorkbook

# Rename Sheet2 to NinzaAutomation
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Name = "NinzaAutomation"

# Clear old data in A1:B1, then move data to F6:G6 (swapped)
$ws.Range("A1").Value = $null
$ws.Range("B1").Value = $null

$ws.Range("F6").Value = "Price"
$ws.Range("G6").Value = "Product Name"
